$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: match the existing bold/bordered/centered header style (same as H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: I and J values per row
$data = @{
    2  = @(4, 5)
    3  = @(7, 8)
    4  = @(8, 8)
    5  = @(8, 8)
    6  = @(7, 7)
    7  = @(6, 6)
    8  = @(7, 7)
    9  = @(6, 7)
    10 = @(8, 8)
    11 = @(7, 8)
    12 = @(8, 8)
    13 = @(6, 6)
    14 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
